$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2469.2222
$ws.Range("I70").Value = 1514.4546
$ws.Range("K70").Value = 4543.3638
$ws.Range("M70").Value = -4273.3638
$ws.Range("H73").Value = 2469.2222
$ws.Range("I73").Value = 1514.4546
$ws.Range("K73").Value = 4543.3638
$ws.Range("M73").Value = -3607.3638
$ws.Range("H92").Value = 870.9545000000001
$ws.Range("J92").Value = 1242.7142
$ws.Range("L92").Value = 1242.7142
$ws.Range("N92").Value = -3738.7142
$ws.Range("H111").Value = 3040
$ws.Range("I111").Value = 2216.1428
$ws.Range("K111").Value = 6648.428400000001
$ws.Range("M111").Value = -3581.428400000001
$ws.Range("H125").Value = 2938.2222
$ws.Range("I125").Value = 940.8333
$ws.Range("J125").Value = 6933
$ws.Range("K125").Value = 8467.4997
$ws.Range("L125").Value = 62397
$ws.Range("M125").Value = -6007.4997
$ws.Range("N125").Value = -67317
$ws.Range("H132").Value = 12323.245
$ws.Range("I132").Value = 2002.7954
$ws.Range("K132").Value = 6008.3862
$ws.Range("M132").Value = -3478.3862
$ws.Range("H138").Value = 4209.4575
$ws.Range("I138").Value = 2295.4285
$ws.Range("J138").Value = 4467.115
$ws.Range("K138").Value = 6886.2855
$ws.Range("L138").Value = 13401.345
$ws.Range("M138").Value = -1746.2855
$ws.Range("N138").Value = -23681.345
$ws.Range("H141").Value = 5583.478
$ws.Range("J141").Value = 29995
$ws.Range("L141").Value = 89985
$ws.Range("N141").Value = -100345

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1231.5
$ws.Range("I32").Value = 1227.7778
$ws.Range("J32").Value = 1600
$ws.Range("K32").Value = 1227.7778
$ws.Range("L32").Value = 1600
$ws.Range("M32").Value = -940.7778000000001
$ws.Range("N32").Value = -2174
$ws.Range("H110").Value = 2533
$ws.Range("I110").Value = 1299.5
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 1299.5
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 745.5
$ws.Range("N110").Value = -9090
$ws.Range("H135").Value = 70000
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1263.32
$ws.Range("I20").Value = 1319.3125
$ws.Range("J20").Value = 1163.7778
$ws.Range("K20").Value = 1319.3125
$ws.Range("L20").Value = 1163.7778
$ws.Range("M20").Value = -1072.3125
$ws.Range("N20").Value = -1657.7778
$ws.Range("H86").Value = 3736.5386
$ws.Range("I86").Value = 3134
$ws.Range("J86").Value = 5270.273
$ws.Range("K86").Value = 3134
$ws.Range("L86").Value = 5270.273
$ws.Range("M86").Value = -2011
$ws.Range("N86").Value = -7516.273
$ws.Range("H89").Value = 3736.5386
$ws.Range("I89").Value = 3134
$ws.Range("J89").Value = 5270.273
$ws.Range("K89").Value = 15670
$ws.Range("L89").Value = 26351.365
$ws.Range("M89").Value = -10054
$ws.Range("N89").Value = -37583.36500000001
$ws.Range("H94").Value = 2120.889
$ws.Range("I94").Value = 2207.8845
$ws.Range("J94").Value = 1894.7
$ws.Range("K94").Value = 2207.8845
$ws.Range("L94").Value = 1894.7
$ws.Range("M94").Value = -1756.8845
$ws.Range("N94").Value = -2796.7
$ws.Range("H99").Value = 30364.5
$ws.Range("I99").Value = 41540.3
$ws.Range("K99").Value = 41540.3
$ws.Range("M99").Value = -40042.3
$ws.Range("H107").Value = 11212.944
$ws.Range("I107").Value = 8455.532999999999
$ws.Range("K107").Value = 8455.532999999999
$ws.Range("M107").Value = -6535.532999999999
$ws.Range("H125").Value = 48454.453
$ws.Range("J125").Value = 48454.453
$ws.Range("L125").Value = 48454.453
$ws.Range("N125").Value = -58294.453
$ws.Range("H134").Value = 2245.05
$ws.Range("I134").Value = 1630.3636
$ws.Range("J134").Value = 5142.857
$ws.Range("K134").Value = 4891.0908
$ws.Range("L134").Value = 15428.571
$ws.Range("M134").Value = -2356.0908
$ws.Range("N134").Value = -20498.571

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 280.92307
$ws.Range("I7").Value = 95
$ws.Range("J7").Value = 314.72726
$ws.Range("K7").Value = 95
$ws.Range("L7").Value = 314.72726
$ws.Range("M7").Value = 18
$ws.Range("N7").Value = -540.72726
$ws.Range("H22").Value = 854.1852
$ws.Range("I22").Value = 477.36365
$ws.Range("J22").Value = 1113.25
$ws.Range("K22").Value = 477.36365
$ws.Range("L22").Value = 1113.25
$ws.Range("M22").Value = -127.36365
$ws.Range("N22").Value = -1813.25
$ws.Range("H109").Value = 32321.25
$ws.Range("J109").Value = 32321.25
$ws.Range("L109").Value = 32321.25
$ws.Range("N109").Value = -34401.25
$ws.Range("H132").Value = 1759.96
$ws.Range("I132").Value = 1791.625
$ws.Range("K132").Value = 5374.875
$ws.Range("M132").Value = -2844.875
$ws.Range("H138").Value = 65149.066
$ws.Range("I138").Value = 55647.4
$ws.Range("J138").Value = 69899.89999999999
$ws.Range("K138").Value = 55647.4
$ws.Range("L138").Value = 69899.89999999999
$ws.Range("M138").Value = -50507.4
$ws.Range("N138").Value = -80179.89999999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1414.8334
$ws.Range("J23").Value = 172.75
$ws.Range("L23").Value = 518.25
$ws.Range("N23").Value = -988.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 62568.58
$ws.Range("I70").Value = 83233.5
$ws.Range("K70").Value = 83233.5
$ws.Range("M70").Value = -82963.5
$ws.Range("H73").Value = 62568.58
$ws.Range("I73").Value = 83233.5
$ws.Range("K73").Value = 83233.5
$ws.Range("M73").Value = -82297.5
$ws.Range("H126").Value = 3145.6667
$ws.Range("I126").Value = 3186.182
$ws.Range("K126").Value = 9558.545999999998
$ws.Range("M126").Value = -7088.545999999998
$ws.Range("H132").Value = 1280.1852
$ws.Range("I132").Value = 1204.68
$ws.Range("K132").Value = 3614.04
$ws.Range("M132").Value = -1084.04

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 523.6667
$ws.Range("I16").Value = 627.75
$ws.Range("K16").Value = 627.75
$ws.Range("M16").Value = -457.75
$ws.Range("H40").Value = 5533.923
$ws.Range("I40").Value = 5694.6523
$ws.Range("K40").Value = 5694.6523
$ws.Range("M40").Value = -5558.6523
$ws.Range("H46").Value = 1337.25
$ws.Range("I46").Value = 760.6
$ws.Range("J46").Value = 2298.3333
$ws.Range("K46").Value = 760.6
$ws.Range("L46").Value = 2298.3333
$ws.Range("M46").Value = -572.6
$ws.Range("N46").Value = -2674.3333
$ws.Range("H68").Value = 3843.0833
$ws.Range("I68").Value = 3786.3333
$ws.Range("K68").Value = 3786.3333
$ws.Range("M68").Value = -3037.3333
$ws.Range("H71").Value = 3843.0833
$ws.Range("I71").Value = 3786.3333
$ws.Range("K71").Value = 18931.6665
$ws.Range("M71").Value = -15187.6665
$ws.Range("H100").Value = 58799.523
$ws.Range("I100").Value = 104645.91
$ws.Range("K100").Value = 104645.91
$ws.Range("M100").Value = -104104.91
$ws.Range("H122").Value = 4298
$ws.Range("I122").Value = 1958.9375
$ws.Range("J122").Value = 6267.737
$ws.Range("K122").Value = 5876.8125
$ws.Range("L122").Value = 18803.211
$ws.Range("M122").Value = -3426.8125
$ws.Range("N122").Value = -23703.211

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3870.75
$ws.Range("I62").Value = 3828
$ws.Range("K62").Value = 3828
$ws.Range("M62").Value = -3204
$ws.Range("H65").Value = 3870.75
$ws.Range("I65").Value = 3828
$ws.Range("K65").Value = 19140
$ws.Range("M65").Value = -16020
$ws.Range("H109").Value = 44761.906
$ws.Range("J109").Value = 44761.906
$ws.Range("L109").Value = 44761.906
$ws.Range("N109").Value = -47535.906
$ws.Range("H122").Value = 2407.037
$ws.Range("I122").Value = 2366.0417
$ws.Range("K122").Value = 7098.125100000001
$ws.Range("M122").Value = -4648.125100000001
